# Append: 2026-02-04 06:58 JST
# Update the "取得日時" (retrieved-at) timestamp stamped on every existing
# data row of the "ランサーズ" sheet from 06:47:46 to 06:58:23 (new scrape run).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2026-02-04 06:47:46"
$newTimestamp = "2026-02-04 06:58:23"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
